$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7272
$ws.Range("C3").Value = 1228
$ws.Range("C4").Value = 918
$ws.Range("C5").Value = 441
$ws.Range("C6").Value = 441
$ws.Range("C7").Value = 436
$ws.Range("C8").Value = 436
$ws.Range("C9").Value = 279
$ws.Range("C10").Value = 279

$ws.Range("B11").Value = "Textiles & Cozy Items"
$ws.Range("C11").Value = 278
